$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("C2")
Write-Host "Before:" $r.Characters().Text
$r.Characters(1,9).Text = "std.error"
Write-Host "After:" $r.Characters().Text
